# Refresh the cryptos price/volume snapshot (scheduled scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "61.591.07" which Excel would otherwise
# reinterpret as a number (dropping the text formatting of the scraped
# value), so keep it formatted as Text before writing the refreshed values.
$ws.Range("D2:D51").NumberFormat = "@"

# New Price (column D) / Volume(1h) (column E) values pulled by the scraper.
$updates = @{
    2 = @{ D = "61.717.65"; E = "  +2.14%  " }
    3 = @{ D = "2.389.59"; E = "  +2.34%  " }
    4 = @{ E = "  -0.03%  " }
    5 = @{ D = "552.09"; E = "  +2.07%  " }
    6 = @{ D = "142.05"; E = "  +4.34%  " }
    7 = @{ E = "  -0.07%  " }
    8 = @{ D = "0.523"; E = "  +0.03%  " }
    9 = @{ D = "2.388.50"; E = "  +2.29%  " }
    10 = @{ E = "  +3.79%  " }
    11 = @{ E = "  +1.90%  " }
    12 = @{ E = "  +1.41%  " }
    13 = @{ E = "  +3.75%  " }
    14 = @{ D = "25.93"; E = "  +6.21%  " }
    15 = @{ E = "  +9.12%  " }
    16 = @{ D = "2.819.24"; E = "  +2.29%  " }
    17 = @{ D = "61.619.92"; E = "  +1.90%  " }
    18 = @{ D = "2.388.09"; E = "  +2.27%  " }
    19 = @{ D = "11.11"; E = "  +5.90%  " }
    20 = @{ D = "322.78"; E = "  +2.01%  " }
    21 = @{ E = "  +2.74%  " }
    22 = @{ D = "6.65"; E = "  +1.76%  " }
    23 = @{ E = "  +0.27%  " }
    24 = @{ D = "64.48"; E = "  +2.51%  " }
    25 = @{ E = "  -5.35%  " }
    26 = @{ D = "9.19"; E = "  +7.23%  " }
    27 = @{ D = "554.75"; E = "  +11.30%  " }
    28 = @{ E = "  +0.06%  " }
    29 = @{ D = "2.500.64"; E = "  +2.06%  " }
    30 = @{ D = "8.28"; E = "  +4.59%  " }
    31 = @{ D = "0.0₃0918"; E = "  +3.39%  " }
    32 = @{ D = "1.42"; E = "  +2.48%  " }
    33 = @{ E = "  +2.81%  " }
    34 = @{ E = "  +3.83%  " }
    35 = @{ D = "1.53"; E = "  +0.80%  " }
    36 = @{ D = "5.75"; E = "  +10.42%  " }
    37 = @{ D = "0.998"; E = "  -0.10%  " }
    38 = @{ D = "1.95"; E = "  +8.81%  " }
    39 = @{ D = "4.74"; E = "  +3.32%  " }
    40 = @{ D = "0.382"; E = "  +2.74%  " }
    41 = @{ D = "18.59"; E = "  +1.97%  " }
    43 = @{ E = "  +0.15%  " }
    44 = @{ D = "148.72"; E = "  +5.57%  " }
    45 = @{ D = "2.26"; E = "  +7.75%  " }
    46 = @{ D = "3.63"; E = "  +3.02%  " }
    47 = @{ D = "0.0529"; E = "  +3.97%  " }
    48 = @{ D = "20.19"; E = "  +4.43%  " }
    49 = @{ D = "0.585"; E = "  +3.36%  " }
    50 = @{ D = "0.0905"; E = "  +1.03%  " }
    51 = @{ D = "0.0224"; E = "  +1.65%  " }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    if ($vals.ContainsKey("D")) {
        $ws.Cells.Item($rowNum, 4).Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($rowNum, 5).Value = $vals["E"]
    }
}
